$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# New rows of asset data (rows 5-7)
# Cell fill order matters for the shared-strings table ordering, so match
# the original author's entry order (B5 before A5, etc).
$ws.Range("B5").Value = "wav"
$ws.Range("A5").Value = "mixkit-game-show-happy-timer-666"
$ws.Range("C5").Value = "Background Game show timer"
$ws.Range("D5").Value = "https://mixkit.co/free-sound-effects/game-show/"

$ws.Range("A6").Value = "mixkit-game-show-suspense-waiting-667"
$ws.Range("B6").Value = "wav"
$ws.Range("C6").Value = "Background Game show timer"
$ws.Range("D6").Value = "https://mixkit.co/free-sound-effects/game-show/"

$ws.Range("A7").Value = "Good-Morning-Doctor-Weird"
$ws.Range("B7").Value = "mp3"
$ws.Range("C7").Value = "Looping lobby/intro music"
$ws.Range("D7").Value = "http://soundimage.org/wp-content/uploads/2016/06/Good-Morning-Doctor-Weird.mp3"

# Column widths (A=56.36328125, B=13.453125, C=32.81640625 target char widths)
$ws.Columns.Item(1).ColumnWidth = 55.5
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 32.0

# Selection after edits
$ws.Range("A8").Select()

# Window geometry (position/size of the workbook window)
$win = $excel.ActiveWindow
$win.Left = 340
$win.Top = 1410
$win.Width = 14400
$win.Height = 7360
